$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("02_38커뮤니케이션(최근일자기준)")

$data = @(
  @("하이젠알앤엠", "2024.06.07~06.13", "4,500~5,500", "-", "15300", "한국투자증권"),
  @("한국스팩14호", "2024.06.03~06.04", "2,000~2,000", "-", "8000", "한국투자증권"),
  @("이노그리드", "2024.05.31~06.07", "29,000~35,000", "-", "17400", "한국투자증권"),
  @("에스오에스랩", "2024.05.28~06.03", "7,500~9,000", "-", "15000", "한국투자증권"),
  @("디비금융스팩12호", "2024.05.28~05.29", "2,000~2,000", "-", "10000", "DB금융투자"),
  @("씨어스테크놀로지", "2024.05.27~05.31", "10,500~14,000", "-", "13650", "한국투자증권"),
  @("라메디텍", "2024.05.27~05.31", "10,400~12,700", "-", "13499", "대신증권"),
  @("미래에셋비전스팩5호", "2024.05.23~05.24", "2,000~2,000", "-", "9500", "미래에셋증권"),
  @("그리드위즈", "2024.05.23~05.29", "34,000~40,000", "-", "47600", "삼성증권"),
  @("이노스페이스", "2024.05.23~05.29", "36,400~45,600", "-", "48412", "미래에셋증권,신한투자증권"),
  @("에이치브이엠(구.한국진공야금)", "2024.05.22~05.28", "11,000~14,200", "-", "26400", "NH투자증권"),
  @("하스", "2024.05.16~05.22", "9,000~12,000", "-", "16290", "삼성증권"),
  @("미래에셋비전스팩4호", "2024.05.13~05.14", "2,000~2,000", "-", "-", "미래에셋증권"),
  @("노브랜드", "2024.04.30~05.08", "8,700~11,000", "-", "10440", "삼성증권"),
  @("KB스팩28호", "2024.04.29~04.30", "2,000~2,000", "2000", "10000", "KB증권"),
  @("아이씨티케이", "2024.04.24~04.30", "13,000~16,000", "20000", "25610", "NH투자증권"),
  @("SK증권스팩12호", "2024.04.17~04.18", "2,000~2,000", "2000", "6000", "SK증권"),
  @("HD현대마린솔루션(구.HD현대글로벌서비스)(유가)", "2024.04.16~04.22", "73,300~83,400", "83400", "652370", "KB증권,신한투자증권,하나증권,대신증권,삼성증권"),
  @("코칩", "2024.04.15~04.19", "11,000~14,000", "18000", "16500", "한국투자증권"),
  @("유안타스팩16호", "2024.04.15~04.16", "2,000~2,000", "2000", "10300", "유안타증권")
)

for ($i = 0; $i -lt $data.Count; $i++) {
  $r = $i + 2
  $row = $data[$i]
  for ($j = 0; $j -lt $row.Count; $j++) {
    $c = $j + 1
    $cell = $ws.Cells.Item($r, $c)
    $cell.Value = "'" + $row[$j]
    $cell.Style = "Normal"
  }
}

Write-Host "done"